# Cosmetics added to Statistics form
# Appends 9 new transaction rows (rows 13-21) to the "Transactions" sheet,
# matching columns: A=Id, B=Amount, C=Date, D=Type, E=Description.
#
# All cells on this sheet are stored as literal text (even the numeric- and
# date-looking ones), so every value is written with a leading apostrophe to
# stop Excel from auto-converting it to a Number/Date, and the cell style is
# then reset back to "Normal" so the quote-prefix formatting doesn't leave a
# stray style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("V4XKUQVU51", "123.484", "08/04/2015", "RegularIncome",   "12455"),
    @("51E139G34B", "233.34",  "08/04/2015", "IrregularExpense","Money"),
    @("2KPOVENKPY", "1239.4",  "08/04/2015", "IrregularExpense","Train"),
    @("IPGVVMHLF9", "12333.4", "08/04/2015", "IrregularExpense","Vacation"),
    @("VTY297HZCK", "234.55",  "08/04/2015", "IrregularExpense","Hotel"),
    @("7FLBKM7TVW", "12.55",   "08/04/2015", "RegularExpense",  "Food"),
    @("WXZWIVXFTG", "123",     "08/04/2015", "IrregularExpense","444"),
    @("1HLZF3HSL1", "123",     "08/04/2015", "RegularIncome",   "123"),
    @("RHTMK6BYYV", "123",     "08/04/2015", "RegularIncome",   "123")
)

$startRow = 13
$columns = @("A", "B", "C", "D", "E")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNumber = $startRow + $i
    $rowValues = $newRows[$i]

    for ($c = 0; $c -lt $columns.Count; $c++) {
        $cell = $ws.Range("$($columns[$c])$rowNumber")
        $cell.Value = "'" + $rowValues[$c]
        $cell.Style = "Normal"
    }
}
